$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 518.0813332500001
$schedule.Range("F2").Value = 34.26463844246032

$schedule.Range("E3").Value = -242.278062
$schedule.Range("F3").Value = -8.011840674603175

$schedule.Range("E4").Value = 489.3200715
$schedule.Range("F4").Value = 32.36243859126984

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B5").Value = 65
$detailed.Range("B6").Value = 64.53478
$detailed.Range("B7").Value = 65.69616000000001
$detailed.Range("C7").Value = "historical"
$detailed.Range("B8").Value = 68.26438
$detailed.Range("C8").Value = "historical"
$detailed.Range("B9").Value = 73.20005
$detailed.Range("C9").Value = "historical"
$detailed.Range("B11").Value = 77.94
$detailed.Range("B12").Value = 79.95026
$detailed.Range("B13").Value = 97.45505
$detailed.Range("B14").Value = 102.9893
$detailed.Range("B15").Value = 78
$detailed.Range("B16").Value = 40.54
$detailed.Range("B17").Value = 5.25447
$detailed.Range("B18").Value = 0.00949
$detailed.Range("B19").Value = -5.62893
$detailed.Range("B20").Value = -6.47747
$detailed.Range("B21").Value = -7.56939
$detailed.Range("B22").Value = -8.945880000000001
$detailed.Range("B23").Value = -10
$detailed.Range("B24").Value = -15.62061
$detailed.Range("B25").Value = -15.6729
$detailed.Range("B26").Value = -16.16569
$detailed.Range("B27").Value = -21.92619
$detailed.Range("B28").Value = -22.22025
$detailed.Range("B29").Value = -24.64712
$detailed.Range("B30").Value = -26.06812
$detailed.Range("B33").Value = -20.55726
$detailed.Range("B37").Value = 48.11341
$detailed.Range("B38").Value = 47.64295
$detailed.Range("B41").Value = 78
$detailed.Range("B44").Value = 57.36105
$detailed.Range("B45").Value = 59.71187
$detailed.Range("B47").Value = 61.93568
$detailed.Range("B48").Value = 63.53195
$detailed.Range("B49").Value = 63.94619
